$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Font-size reductions on the three big centered title paragraphs
#    (both "sz" and complex-script "szCs" need to match; Font.Size
#    drives w:sz, Font.SizeBi drives w:szCs in this host) and the
#    paragraph-mark run-properties (w:pPr/w:rPr) pick up the same
#    explicit size.
# ------------------------------------------------------------------

# "***=licencia.tipo_lic*** de uso e aproveitamento de Água ***=licencia.tipo_agua***"
# 58 half-points -> 52 half-points (26pt)
$pTipo = $d.Paragraphs(8)
$pTipo.Range.Font.Size = 26
$pTipo.Range.Font.SizeBi = 26

# "***=licencia.lic_nro***"
# sz 52 / szCs 60 -> 48 half-points (24pt) both
$pLicNro = $d.Paragraphs(10)
$pLicNro.Range.Font.Size = 24
$pLicNro.Range.Font.SizeBi = 24

# "***=exp_name***"
# sz 52 / szCs 60 -> 40 half-points (20pt) both
$pExpName = $d.Paragraphs(11)
$pExpName.Range.Font.Size = 20
$pExpName.Range.Font.SizeBi = 20

# ------------------------------------------------------------------
# 2. Relocate the lone "_GoBack" bookmark from the "OBSERVAÇÕES" title
#    (4th section header) to the very end of the exp_name paragraph
#    (right after its last run, before the paragraph mark). Word only
#    ever keeps a single "_GoBack" bookmark, so adding it at the new
#    spot implicitly removes the old one.
#
#    A zero-length Range placed exactly at the end-of-paragraph offset
#    is mis-resolved by this host, so a tiny marker run is inserted at
#    that exact spot first (InsertAfter on a collapsed Find range does
#    land correctly), the bookmark is anchored to the marker's Range,
#    and the marker text is then cleared back out - leaving a clean
#    zero-length bookmark exactly where it belongs.
# ------------------------------------------------------------------
$markerTag = "ZZ_GOBACK_MARKER_ZZ"

$findRng = $d.Content
$findRng.Find.Execute("***=exp_name***", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.Collapse(0)
$findRng.InsertAfter($markerTag)

$markerRng = $d.Content
$markerRng.Find.Execute($markerTag, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRng) | Out-Null

$clearRng = $d.Content
$clearRng.Find.Execute($markerTag, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$clearRng.Text = ""

# ------------------------------------------------------------------
# 3. Merge the two runs of the "4    OBSERVAÇÕES:" heading (previously
#    split apart by the bookmark that just moved away) back into a
#    single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("4    OBSE" + "RVAÇÕES:", $true, $false, $false, $false, $false, $true, 1, $false, `
    "4    OBSERVAÇÕES:", 2) | Out-Null
